$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in T2 (131518 -> 136312)
$ws.Range("T2").Value = 136312

# Move the active selection from T3 to T2
$ws.Range("T2").Select()
